$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.382.49"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.841.63"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'239.29"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.6262"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.07427"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "'0.2892"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'24.94"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").Value = "'0.07720"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.843.91"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'4.974"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'0.6741"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'6.212"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "29.391.37"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'233.93"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").Value = "'12.30"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'7.287"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'158.45"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'8.483"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "'0.1345"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "'0.07253"
$ws.Range("E28").Value = "  +13.69%  "
$ws.Range("D29").Value = "'1.466"
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("D30").Value = "'1.477"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'4.053"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'4.024"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").Value = "'1.140"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'0.6972"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'2.571"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").Value = "'0.01841"
$ws.Range("D38").Value = "'6.919"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("D39").Value = "'2.817"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "1.232.46"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "'0.9644"
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "'101.08"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'65.44"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'6.955"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "'8.909"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'0.1135"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").Value = "'0.3896"
$ws.Range("E51").Value = "  -1.61%  "